# Reorder the "Recorded By" entries in column G so that the literal
# "System" token is moved to the end of the comma-separated list. If the
# list does not contain the literal "System" token (case-sensitive),
# the whole list order is simply reversed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = 7
    $val = $cell.Value2

    if ($null -eq $val) { continue }
    if ($val -notmatch ",") { continue }

    $parts = $val -split ","
    $trimmed = @()
    foreach ($p in $parts) { $trimmed += $p.Trim() }

    # Case-sensitive check for the literal token "System" (PowerShell's
    # built-in -eq/-ne/-contains operators are case-insensitive in this
    # runtime, so use the .NET string .Equals() instance method, which
    # performs an ordinal, case-sensitive comparison).
    $hasSystem = $false
    foreach ($p in $trimmed) {
        if ($p.Equals("System")) { $hasSystem = $true }
    }

    if ($hasSystem) {
        $rest = @()
        foreach ($p in $trimmed) {
            if (-not $p.Equals("System")) { $rest += $p }
        }
        $newParts = $rest + @("System")
    } else {
        $newParts = $trimmed[($trimmed.Length - 1)..0]
    }

    $newVal = [string]::Join(", ", $newParts)
    $cell.Value2 = $newVal
}
